# reservations.xlsx — "Add files via upload" edit
#
# 1. Sheet1: the "charges" (P) / "%" (Q) commission columns are zeroed out
#    for every reservation row.
# 2. Plateformes sheet: the bold/centred header styling is stripped, and
#    three platform -> colour rows (Booking/Airbnb/Autre) are added below
#    the existing header row.
# 3. The Plateformes sheet becomes the active tab/selection instead of
#    Sheet1.

$wb = $excel.ActiveWorkbook

# --- Sheet1: zero out columns P ("charges") and Q ("%") for all data rows ---
$ws1 = $wb.Worksheets.Item("Sheet1")
$lastRow = $ws1.Cells.Item($ws1.Rows.Count, 1).End(-4162).Row   # xlUp
$ws1.Range("P2:Q" + $lastRow).Value = 0
[void]$ws1.Range("A1").Select()

# --- Plateformes sheet: drop the old bold header font/alignment ---
$ws2 = $wb.Worksheets.Item("Plateformes")
$ws2.Rows.Item(1).ClearFormats()

# add the platform / colour rows
$ws2.Range("A2").Value = "Booking"
$ws2.Range("B2").Value = "#1e90ff"
$ws2.Rows.Item(2).RowHeight = 12.5

$ws2.Range("A3").Value = "Airbnb"
$ws2.Range("B3").Value = "#e74c3c"

$ws2.Range("A4").Value = "Autre"
$ws2.Range("B4").Value = "#f59e0b"

# Plateformes becomes the selected/active sheet, with B4 selected
[void]$ws2.Range("B4").Select()
[void]$ws2.Activate()
